# Add season-record columns (Wins, Losses, Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1:AF1 - copy the formatting of the adjacent header
# cell (AC1, style index 1 - bold/centered/bordered) onto the new range,
# then set the header labels.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-53) shares the same team record for the season.
for ($row = 2; $row -le 53; $row++) {
    $ws.Cells.Item($row, 30).Value = 111
    $ws.Cells.Item($row, 31).Value = 51
    $ws.Cells.Item($row, 32).Value = 0
}
